$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = "64.060.66"
$c.ClearFormats()
$ws.Cells.Item(2,5).Value = "  +0.37%  "
$c = $ws.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = "3.133.01"
$c.ClearFormats()
$ws.Cells.Item(3,5).Value = "  +0.03%  "
$ws.Cells.Item(4,5).Value = "  +0.02%  "
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "587.35"
$c.ClearFormats()
$ws.Cells.Item(5,5).Value = "  -0.60%  "
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "146.18"
$c.ClearFormats()
$ws.Cells.Item(6,5).Value = "  -0.32%  "
$ws.Cells.Item(7,5).Value = "  +0.02%  "
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "3.129.70"
$c.ClearFormats()
$ws.Cells.Item(8,5).Value = "  +0.17%  "
$ws.Cells.Item(9,5).Value = "  -1.49%  "
$ws.Cells.Item(10,5).Value = "  -2.72%  "
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "5.83"
$c.ClearFormats()
$ws.Cells.Item(11,5).Value = "  +2.22%  "
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "0.457"
$c.ClearFormats()
$ws.Cells.Item(12,5).Value = "  -2.17%  "
$ws.Cells.Item(13,5).Value = "  -3.35%  "
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "37.09"
$c.ClearFormats()
$ws.Cells.Item(14,5).Value = "  +2.97%  "
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "3.657.47"
$c.ClearFormats()
$ws.Cells.Item(15,5).Value = "  +0.18%  "
$ws.Cells.Item(16,5).Value = "  -1.51%  "
$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = "63.862.46"
$c.ClearFormats()
$ws.Cells.Item(17,5).Value = "  +0.17%  "
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = "3.131.01"
$c.ClearFormats()
$ws.Cells.Item(18,5).Value = "  +0.07%  "
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "7.08"
$c.ClearFormats()
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "463.08"
$c.ClearFormats()
$ws.Cells.Item(20,5).Value = "  -0.62%  "
$ws.Cells.Item(21,5).Value = "  +0.23%  "
$ws.Cells.Item(22,5).Value = "  -0.59%  "
$ws.Cells.Item(23,5).Value = "  -2.15%  "
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "12.86"
$c.ClearFormats()
$ws.Cells.Item(24,5).Value = "  -3.30%  "
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = "80.69"
$c.ClearFormats()
$ws.Cells.Item(25,5).Value = "  -2.00%  "
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "2.30"
$c.ClearFormats()
$ws.Cells.Item(26,5).Value = "  +6.70%  "
$ws.Cells.Item(27,5).Value = "  -0.01%  "
$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value = "9.50"
$c.ClearFormats()
$ws.Cells.Item(28,5).Value = "  +9.60%  "
$ws.Cells.Item(29,5).Value = "  -1.25%  "
$ws.Cells.Item(30,5).Value = "  +0.14%  "
$ws.Cells.Item(31,5).Value = "  -1.21%  "
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "7.13"
$c.ClearFormats()
$ws.Cells.Item(32,5).Value = "  +3.90%  "
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.ClearFormats()
$ws.Cells.Item(34,5).Value = "  -0.11%  "
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0851"
$c.ClearFormats()
$ws.Cells.Item(35,5).Value = "  -2.38%  "
$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value = "1.05"
$c.ClearFormats()
$ws.Cells.Item(36,5).Value = "  +0.05%  "
$ws.Cells.Item(37,5).Value = "  -3.86%  "
$ws.Cells.Item(38,2).Value = "Filecoin"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = "6.03"
$c.ClearFormats()
$ws.Cells.Item(38,5).Value = "  -1.61%  "
$ws.Cells.Item(39,2).Value = "dogwifhat"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.ClearFormats()
$ws.Cells.Item(39,5).Value = "  -3.21%  "
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "51.35"
$c.ClearFormats()
$ws.Cells.Item(40,5).Value = "  +1.07%  "
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = "437.85"
$c.ClearFormats()
$ws.Cells.Item(41,5).Value = "  -3.05%  "
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "8.89"
$c.ClearFormats()
$ws.Cells.Item(42,5).Value = "  +2.08%  "
$ws.Cells.Item(43,2).Value = "VeChain"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "0.0371"
$c.ClearFormats()
$ws.Cells.Item(43,5).Value = "  -0.73%  "
$ws.Cells.Item(44,2).Value = "Maker"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "2.904.58"
$c.ClearFormats()
$ws.Cells.Item(44,5).Value = "  -0.53%  "
$ws.Cells.Item(45,2).Value = "TheGraph"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = "0.283"
$c.ClearFormats()
$ws.Cells.Item(45,5).Value = "  +1.99%  "
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = "39.63"
$c.ClearFormats()
$ws.Cells.Item(46,5).Value = "  +15.29%  "
$ws.Cells.Item(47,5).Value = "  -3.39%  "
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "126.69"
$c.ClearFormats()
$ws.Cells.Item(48,5).Value = "  -1.24%  "
$ws.Cells.Item(49,5).Value = "  +0.00%  "
$ws.Cells.Item(50,5).Value = "  -0.98%  "
$ws.Cells.Item(51,2).Value = "InjectiveProtocol"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "24.11"
$c.ClearFormats()
$ws.Cells.Item(51,5).Value = "  -2.27%  "
